$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column M ("Gas Choice ID"), shifting
# Gas Choice ID / Gas Rate Code / Gas Usage (therms) one column to the right
# (M->N, N->O, O->P). The sheet's used range grows from A1:O2 to A1:P2.
$ws.Columns("M").Insert()

# New header cell for the inserted "Gas Supplier" column, matching the
# formatting used by the rest of the header row (bold font, thin border,
# centered horizontally, top-aligned vertically).
$ws.Range("M1").Value = "Gas Supplier"
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108
$ws.Range("M1").VerticalAlignment = -4160
$ws.Range("M1").Borders.LineStyle = 1

# New data value for the inserted column on row 2.
$ws.Range("M2").Value = "N/A"
